# Apply narrow (0.5") margins to every section and widen the two
# template tables' columns from 5040 (3.5") to 5400 twips (3.75")
# each -- i.e. 252pt -> 270pt, expressed in points since Word's COM
# object model reports/accepts table and margin measurements in points
# (1 pt = 20 twips).

$d = $word.ActiveDocument

# --- Widen the columns of both tables in the document ---------------
$tables = $d.Tables
for ($i = 1; $i -le $tables.Count; $i++) {
    $tbl = $tables.Item($i)
    $cols = $tbl.Columns
    for ($c = 1; $c -le $cols.Count; $c++) {
        $cols.Item($c).Width = 270
    }
}

# --- Narrow every section's page margins to 0.5" (36pt / 720 twips) -
foreach ($sec in $d.Sections) {
    $sec.PageSetup.TopMargin = 36
    $sec.PageSetup.BottomMargin = 36
    $sec.PageSetup.LeftMargin = 36
    $sec.PageSetup.RightMargin = 36
}
